# Acorto nombres de nodos con demasiados caracteres
# Shorten a handful of overly long node names in the "maestro-nodos" sheet,
# and tidy up the sheet view / column widths that LibreOffice re-wrote when
# the file was last edited.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Shorten the long node names (column B) that exceeded the character limit
# ---------------------------------------------------------------------------
$ws.Range("B4").Value  = "Gas Natural"
$ws.Range("B15").Value = "Gas de Red"
$ws.Range("B21").Value = "Plantas de Gas"
$ws.Range("B23").Value = "Otros Centros"
$ws.Range("B26").Value = "No Energético"

# ---------------------------------------------------------------------------
# 2. Re-scroll the sheet back to the top and move the active selection to B21
# ---------------------------------------------------------------------------
$excel.ActiveWindow.ScrollRow    = 1
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("B21").Select()

# ---------------------------------------------------------------------------
# 3. Adjust column widths: narrow column B, give column A (and the rest of
#    the sheet up to column AMK/1025) the default width
# ---------------------------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 7.671768707482993
$ws.Columns.Item(2).ColumnWidth = 17.202380952380953
$ws.Range("C1:AMK1").EntireColumn.ColumnWidth = 7.671768707482993
